# Replace the 96 forecast rows (A2:B97) with the newly published Entsoe
# consumption forecast data (Horeco added to the portfolio).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Forecasted Consumption (MW)" values for rows 2-97
$consumption = @(
    5610, 5550, 5490, 5450, 5410, 5390, 5370, 5350, 5340, 5330,
    5320, 5310, 5300, 5290, 5290, 5290, 5310, 5330, 5370, 5420,
    5480, 5540, 5610, 5680, 5750, 5820, 5870, 5920, 5960, 5990,
    6000, 6000, 5980, 5940, 5900, 5850, 5790, 5740, 5690, 5640,
    5610, 5580, 5570, 5560, 5560, 5560, 5560, 5560, 5570, 5590,
    5610, 5630, 5660, 5700, 5740, 5770, 5800, 5830, 5870, 5910,
    5950, 6000, 6060, 6120, 6200, 6270, 6330, 6400, 6460, 6530,
    6610, 6680, 6760, 6830, 6900, 6970, 7040, 7090, 7140, 7220,
    7260, 7260, 7230, 7160, 7030, 6930, 6770, 6620, 6470, 6330,
    6200, 6070, 5840, 5810, 5730, 5640
)

# New "Timestamp" (serial date) values for rows 2-97
$timestamps = @(
    45875, 45875.01041666666, 45875.02083333334, 45875.03125, 45875.04166666666, 45875.05208333334, 45875.0625, 45875.07291666666, 45875.08333333334, 45875.09375,
    45875.10416666666, 45875.11458333334, 45875.125, 45875.13541666666, 45875.14583333334, 45875.15625, 45875.16666666666, 45875.17708333334, 45875.1875, 45875.19791666666,
    45875.20833333334, 45875.21875, 45875.22916666666, 45875.23958333334, 45875.25, 45875.26041666666, 45875.27083333334, 45875.28125, 45875.29166666666, 45875.30208333334,
    45875.3125, 45875.32291666666, 45875.33333333334, 45875.34375, 45875.35416666666, 45875.36458333334, 45875.375, 45875.38541666666, 45875.39583333334, 45875.40625,
    45875.41666666666, 45875.42708333334, 45875.4375, 45875.44791666666, 45875.45833333334, 45875.46875, 45875.47916666666, 45875.48958333334, 45875.5, 45875.51041666666,
    45875.52083333334, 45875.53125, 45875.54166666666, 45875.55208333334, 45875.5625, 45875.57291666666, 45875.58333333334, 45875.59375, 45875.60416666666, 45875.61458333334,
    45875.625, 45875.63541666666, 45875.64583333334, 45875.65625, 45875.66666666666, 45875.67708333334, 45875.6875, 45875.69791666666, 45875.70833333334, 45875.71875,
    45875.72916666666, 45875.73958333334, 45875.75, 45875.76041666666, 45875.77083333334, 45875.78125, 45875.79166666666, 45875.80208333334, 45875.8125, 45875.82291666666,
    45875.83333333334, 45875.84375, 45875.85416666666, 45875.86458333334, 45875.875, 45875.88541666666, 45875.89583333334, 45875.90625, 45875.91666666666, 45875.92708333334,
    45875.9375, 45875.94791666666, 45875.95833333334, 45875.96875, 45875.97916666666, 45875.98958333334
)

$arr = New-Object 'object[,]' $consumption.Count,2
for ($i = 0; $i -lt $consumption.Count; $i++) {
    $arr[$i,0] = $consumption[$i]
    $arr[$i,1] = $timestamps[$i]
}

$ws.Range("A2:B97").Value = $arr

Write-Host "Updated rows 2 to 97 with new consumption forecast data"